$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 20 (2025-07) stats
$ws.Range("B20").Value = 6187
$ws.Range("D20").Value = 5579499
$ws.Range("E20").Value = 901.8100856634879
$ws.Range("F20").Value = 6.875107963378824
$ws.Range("H20").Value = 26.213791429495
